# Weekly update for "Hortaliza, Vega Central Mapocho de Santiago - Poroto granado":
# a new daily price record is inserted as row 147 (pushing the existing
# rows 147-231 down to 148-232), extending the used range to A1:R232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 147, shifting rows 147:231 -> 148:232.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new market record.
$ws.Range("A147").Value = 9
$ws.Range("B147").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C147").Value = "Metropolitana"
$ws.Range("D147").Value = 44606
$ws.Range("E147").Value = 13
$ws.Range("F147").Value = 100112030
$ws.Range("G147").Value = "Poroto granado"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 97
$ws.Range("K147").Value = 30000
$ws.Range("L147").Value = 32000
$ws.Range("M147").Value = 30990
$ws.Range("N147").Value = "`$/saco 25 kilos"
$ws.Range("O147").Value = "Región Metropolitana"
$ws.Range("P147").Value = 1240
$ws.Range("Q147").Value = 25
$ws.Range("R147").Value = "Hortaliza"
